{"js": "// Apply the commit's text replacements to the document body.\n// Each entry is [old exact text, new exact text]. The title/description\n// strings appear more than once, so we search document-wide and replace\n// every match found (search results preserve each run's own formatting).\nconst replacements = [\n  [\n    \"Play Blazin Hot 7s Stack Em Up Slot for Free | Review 2021\",\n    \"Play Blazin' Hot 7s Stack'Em Up for Free\",\n  ],\n  [\n    \"Innovative gameplay with 4 stacked grids.\",\n    \"Unique gameplay with stacked grids\",\n  ],\n  [\n    \"Cascading feature for frequent wins and progression.\",\n    \"Progression mechanics with cascading symbols\",\n  ],\n  [\n    \"Stack'Em Up feature activates the upper grid for more wins.\",\n    \"Exciting 'Stack'Em Up' feature\",\n  ],\n  [\n    \"High volatility provides the chance for big payouts.\",\n    \"High volatility and decent RTP rate\",\n  ],\n  [\n    \"Overwhelming interface on smaller screens.\",\n    \"Overwhelming interface\",\n  ],\n  [\n    \"Limited number of paylines may not appeal to all players.\",\n    \"Challenging navigation for small screen devices\",\n  ],\n  [\n    \"Read our unbiased review of Blazin Hot 7s Stack Em Up slot. Learn how to play the game and try it for free. Discover pros and cons and RTP rate.\",\n    \"Read our review of Blazin' Hot 7s Stack'Em Up and play for free. Discover unique gameplay and exciting features.\",\n  ],\n];\n\nconst body = context.document.body;\nconst resultSets = [];\n\nfor (const [oldText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  resultSets.push(results);\n}\n\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, newText] = replacements[i];\n  const results = resultSets[i];\n  for (let j = 0; j < results.items.length; j++) {\n    results.items[j].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply the commit's text replacements to the document body.\n# We locate each exact phrase with Range.Find and overwrite it by assigning\n# directly to Range.Text (NOT Find.Execute's Replace argument) so Word's\n# smart-quote AutoFormat never touches the straight apostrophes we need.\n# Because the title/description string appears twice, we keep searching\n# forward from the end of each match until no more matches are found.\n\nfunction Replace-AllExact {\n    param(\n        [string]$OldText,\n        [string]$NewText\n    )\n\n    $d = $word.ActiveDocument\n    $searchRange = $d.Content.Duplicate\n    $searchRange.Start = 0\n\n    while ($true) {\n        $find = $searchRange.Find\n        $find.ClearFormatting()\n        $find.Text = $OldText\n        $matchFound = $find.Execute()\n        if (-not $matchFound) {\n            break\n        }\n        $searchRange.Text = $NewText\n        $searchRange.Collapse(0)\n        $searchRange.End = $d.Content.End\n    }\n}\n\nReplace-AllExact \"Play Blazin Hot 7s Stack Em Up Slot for Free | Review 2021\" \"Play Blazin' Hot 7s Stack'Em Up for Free\"\nReplace-AllExact \"Innovative gameplay with 4 stacked grids.\" \"Unique gameplay with stacked grids\"\nReplace-AllExact \"Cascading feature for frequent wins and progression.\" \"Progression mechanics with cascading symbols\"\nReplace-AllExact \"Stack'Em Up feature activates the upper grid for more wins.\" \"Exciting 'Stack'Em Up' feature\"\nReplace-AllExact \"High volatility provides the chance for big payouts.\" \"High volatility and decent RTP rate\"\nReplace-AllExact \"Overwhelming interface on smaller screens.\" \"Overwhelming interface\"\nReplace-AllExact \"Limited number of paylines may not appeal to all players.\" \"Challenging navigation for small screen devices\"\nReplace-AllExact \"Read our unbiased review of Blazin Hot 7s Stack Em Up slot. Learn how to play the game and try it for free. Discover pros and cons and RTP rate.\" \"Read our review of Blazin' Hot 7s Stack'Em Up and play for free. Discover unique gameplay and exciting features.\"\n"}
